$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 135
$ws.Range("I2").Value = 369
$ws.Range("J2").Value = 1468
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = 261
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = 154
$ws.Range("T2").Value = 289
$ws.Range("U2").Value = 27
$ws.Range("V2").Value = 2379
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2345
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 10
